$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 16.6019354869764
$ws.Range("D2").Value = 3.467017911033662
$ws.Range("E2").Value = 24.63841689047146
$ws.Range("F2").Value = 20.57066559520825
$ws.Range("G2").Value = 3.588673489366149
$ws.Range("I2").Value = 24.9666799838896
$ws.Range("L2").Value = 9.896296048513905
$ws.Range("M2").Value = 15.62146801757739
$ws.Range("O2").Value = 18.03288912875454
$ws.Range("B3").Value = 16.10415831715776
$ws.Range("D3").Value = 3.445389423039186
$ws.Range("E3").Value = 24.06936471691901
$ws.Range("F3").Value = 20.43445330180801
$ws.Range("G3").Value = 3.591225691608983
$ws.Range("I3").Value = 25.02522303257962
$ws.Range("L3").Value = 9.718442191773867
$ws.Range("M3").Value = 15.37279268382949
$ws.Range("O3").Value = 18.01100373746687
$ws.Range("B4").Value = 15.79031181921266
$ws.Range("D4").Value = 3.432085857482968
$ws.Range("E4").Value = 23.71461400061439
$ws.Range("F4").Value = 20.35946225689755
$ws.Range("G4").Value = 3.59287569357571
$ws.Range("I4").Value = 25.06957804314576
$ws.Range("L4").Value = 9.608000852797932
$ws.Range("M4").Value = 15.21845276225383
$ws.Range("O4").Value = 18.00485329689827
$ws.Range("B5").Value = 15.66051139723302
$ws.Range("D5").Value = 3.426661484244174
$ws.Range("E5").Value = 23.56889865757917
$ws.Range("F5").Value = 20.33110557422922
$ws.Range("G5").Value = 3.593569008129709
$ws.Range("I5").Value = 25.08976389421926
$ws.Range("L5").Value = 9.562732240965998
$ws.Range("M5").Value = 15.15520354197395
$ws.Range("O5").Value = 18.00417668586953
$ws.Range("B6").Value = 15.63884794735615
$ws.Range("D6").Value = 3.425760669705674
$ws.Range("E6").Value = 23.54463903195578
$ws.Range("F6").Value = 20.32653076639948
$ws.Range("G6").Value = 3.593685398368152
$ws.Range("I6").Value = 25.09324310141585
$ws.Range("L6").Value = 9.555200989851151
$ws.Range("M6").Value = 15.14468145822749
$ws.Range("O6").Value = 18.00417472340259
$ws.Range("B7").Value = 15.78856878102499
$ws.Range("D7").Value = 3.432012711125546
$ws.Range("E7").Value = 23.71265323212858
$ws.Range("F7").Value = 20.35907087497158
$ws.Range("G7").Value = 3.592884959040727
$ws.Range("I7").Value = 25.06984173558011
$ws.Range("L7").Value = 9.60739134374748
$ws.Range("M7").Value = 15.2176011147485
$ws.Range("O7").Value = 18.0048367686466
$ws.Range("B8").Value = 16.43209709355072
$ws.Range("D8").Value = 3.459567027874412
$ws.Range("E8").Value = 24.44343153042697
$ws.Range("F8").Value = 20.52192461200129
$ws.Range("G8").Value = 3.58953631655526
$ws.Range("I8").Value = 24.9851182466123
$ws.Range("L8").Value = 9.835257735287826
$ws.Range("M8").Value = 15.53610386579584
$ws.Range("O8").Value = 18.02382963948289
$ws.Range("B9").Value = 17.62259153534538
$ws.Range("D9").Value = 3.513288019715081
$ws.Range("E9").Value = 25.82601405973585
$ws.Range("F9").Value = 20.90834086720858
$ws.Range("G9").Value = 3.58362455237696
$ws.Range("I9").Value = 24.88583304773672
$ws.Range("L9").Value = 10.2702221295821
$ws.Range("M9").Value = 16.14497427442528
$ws.Range("O9").Value = 18.1189082470885
$ws.Range("B10").Value = 18.44611257498727
$ws.Range("D10").Value = 3.552398220379946
$ws.Range("E10").Value = 26.80104031113454
$ws.Range("F10").Value = 21.23075955409331
$ws.Range("G10").Value = 3.579675988737898
$ws.Range("I10").Value = 24.85377375559928
$ws.Range("L10").Value = 10.57991974365259
$ws.Range("M10").Value = 16.57937750623676
$ws.Range("O10").Value = 18.22387989354285
$ws.Range("B11").Value = 18.80829651811197
$ws.Range("D11").Value = 3.570076973713451
$ws.Range("E11").Value = 27.23386635661063
$ws.Range("F11").Value = 21.38521267722587
$ws.Range("G11").Value = 3.577964467457368
$ws.Range("I11").Value = 24.84807448048748
$ws.Range("L11").Value = 10.71814911798584
$ws.Range("M11").Value = 16.7735218242447
$ws.Range("O11").Value = 18.27917977566029
$ws.Range("B12").Value = 18.94355987704277
$ws.Range("D12").Value = 3.576752246124429
$ws.Range("E12").Value = 27.39608311991638
$ws.Range("F12").Value = 21.44476512841831
$ws.Range("G12").Value = 3.577328466590015
$ws.Range("I12").Value = 24.84719278066573
$ws.Range("L12").Value = 10.77007118290776
$ws.Range("M12").Value = 16.84648785897968
$ws.Range("O12").Value = 18.30119530783939
$ws.Range("B13").Value = 18.91451391650973
$ws.Range("D13").Value = 3.575315513475266
$ws.Range("E13").Value = 27.36122394839333
$ws.Range("F13").Value = 21.43189300422257
$ws.Range("G13").Value = 3.577464902959844
$ws.Range("I13").Value = 24.84732592722193
$ws.Range("L13").Value = 10.75890824837761
$ws.Range("M13").Value = 16.83079870151249
$ws.Range("O13").Value = 18.29640626778269
$ws.Range("B14").Value = 18.81946307987323
$ws.Range("D14").Value = 3.570626556537842
$ws.Range("E14").Value = 27.24724658126629
$ws.Range("F14").Value = 21.39009108456545
$ws.Range("G14").Value = 3.577911900847634
$ws.Range("I14").Value = 24.84797636986942
$ws.Range("L14").Value = 10.72242946844575
$ws.Range("M14").Value = 16.77953614806455
$ws.Range("O14").Value = 18.28096953828096
$ws.Range("B15").Value = 18.76099309197123
$ws.Range("D15").Value = 3.567751832006842
$ws.Range("E15").Value = 27.17720863673543
$ws.Range("F15").Value = 21.36462314530003
$ws.Range("G15").Value = 3.578187275791908
$ws.Range("I15").Value = 24.84854097053139
$ws.Range("L15").Value = 10.70002898238425
$ws.Range("M15").Value = 16.74806289686732
$ws.Range("O15").Value = 18.2716536922388
$ws.Range("B16").Value = 18.42218396831776
$ws.Range("D16").Value = 3.551240340729446
$ws.Range("E16").Value = 26.77252587655841
$ws.Range("F16").Value = 21.22081767627184
$ws.Range("G16").Value = 3.57978953952263
$ws.Range("I16").Value = 24.8543249279855
$ws.Range("L16").Value = 10.57082918926931
$ws.Range("M16").Value = 16.56661530026283
$ws.Range("O16").Value = 18.2204169355905
$ws.Range("B17").Value = 18.21107546496149
$ws.Range("D17").Value = 3.541080015173833
$ws.Range("E17").Value = 26.5214149219175
$ws.Range("F17").Value = 21.13455306236944
$ws.Range("G17").Value = 3.580794124005005
$ws.Range("I17").Value = 24.86014831922439
$ws.Range("L17").Value = 10.49085986802866
$ws.Range("M17").Value = 16.45437595576641
$ws.Range("O17").Value = 18.19091114132976
$ws.Range("B18").Value = 18.08848573462144
$ws.Range("D18").Value = 3.535225642669552
$ws.Range("E18").Value = 26.37598154472654
$ws.Range("F18").Value = 21.08567201069607
$ws.Range("G18").Value = 3.581379910642515
$ws.Range("I18").Value = 24.8643341791642
$ws.Range("L18").Value = 10.44461615254343
$ws.Range("M18").Value = 16.3894954780524
$ws.Range("O18").Value = 18.17465140931549
$ws.Range("B19").Value = 18.04678194344751
$ws.Range("D19").Value = 3.533241757134594
$ws.Range("E19").Value = 26.32657295042302
$ws.Range("F19").Value = 21.06924977148116
$ws.Range("G19").Value = 3.581579619833784
$ws.Range("I19").Value = 24.8658951082847
$ws.Range("L19").Value = 10.42891760237756
$ws.Range("M19").Value = 16.36747419682224
$ws.Range("O19").Value = 18.169268583802
$ws.Range("B20").Value = 18.23366972294822
$ws.Range("D20").Value = 3.542162700220777
$ws.Range("E20").Value = 26.54825078943037
$ws.Range("F20").Value = 21.14366028522222
$ws.Range("G20").Value = 3.580686359232969
$ws.Range("I20").Value = 24.85944184973661
$ws.Range("L20").Value = 10.49939865476401
$ws.Range("M20").Value = 16.46635788975788
$ws.Range("O20").Value = 18.19397854486008
$ws.Range("B21").Value = 18.84743378554712
$ws.Range("D21").Value = 3.572004364788392
$ws.Range("E21").Value = 27.28077131788828
$ws.Range("F21").Value = 21.40234088428611
$ws.Range("G21").Value = 3.577780278445693
$ws.Range("I21").Value = 24.84775068961532
$ws.Range("L21").Value = 10.73315594777086
$ws.Range("M21").Value = 16.79460860281639
$ws.Range("O21").Value = 18.2854746104472
$ws.Range("B22").Value = 19.2375218490617
$ws.Range("D22").Value = 3.591393780604177
$ws.Range("E22").Value = 27.74963534859696
$ws.Range("F22").Value = 21.57758054989772
$ws.Range("G22").Value = 3.575951573003269
$ws.Range("I22").Value = 24.8475489405075
$ws.Range("L22").Value = 10.88345052058663
$ws.Range("M22").Value = 17.00589779747967
$ws.Range("O22").Value = 18.35152998431456
$ws.Range("B23").Value = 19.03036482295195
$ws.Range("D23").Value = 3.581056713469481
$ws.Range("E23").Value = 27.50034239475283
$ws.Range("F23").Value = 21.48350543697156
$ws.Range("G23").Value = 3.576921150118744
$ws.Range("I23").Value = 24.84697658324456
$ws.Range("L23").Value = 10.80347521185605
$ws.Range("M23").Value = 16.89344238404132
$ws.Range("O23").Value = 18.31570656121698
$ws.Range("B24").Value = 18.22345865327844
$ws.Range("D24").Value = 3.541673258752723
$ws.Range("E24").Value = 26.53612160415645
$ws.Range("F24").Value = 21.13954068269129
$ws.Range("G24").Value = 3.580735054021277
$ws.Range("I24").Value = 24.85975863444748
$ws.Range("L24").Value = 10.49553910179103
$ws.Range("M24").Value = 16.4609419514893
$ws.Range("O24").Value = 18.19258957939146
$ws.Range("B25").Value = 17.30901201243531
$ws.Range("D25").Value = 3.498805887710879
$ws.Range("E25").Value = 25.45846678685344
$ws.Range("F25").Value = 20.79685311673293
$ws.Range("G25").Value = 3.585154189114408
$ws.Range("I25").Value = 24.90551783984371
$ws.Range("L25").Value = 10.15410349878396
$ws.Range("M25").Value = 15.98229517787934
$ws.Range("O25").Value = 18.08700041833216
